# Ising_Results.xlsx -- "Added some Data Analysis"
#
# Adds four new columns (G:J) to Sheet1 containing Speedup / Efficiency
# analysis derived from the existing Running Time (column E) data:
#   G = Speedup (1:4)                 = E(1 thread) / E(4 threads)
#   H = Speedup (1:8)                 = E(1 thread) / E(8 threads)
#   I = Efficency (Speed: 4 Thread)   = G / 4
#   J = Efficency (Speed: 8 Thread)   = H / 8
#
# Rows with B=256 have no 8-thread counterpart in the source data, so
# those rows get "N/A" for both Speedup columns instead of a formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each "File" block spans 11 rows:
#   block+0  : B=256,   F=1  -> no 4/8-thread comparison exists (B=256 @ F=8
#                               is missing from the source data)
#   block+1  : B=1024,  F=1
#   block+2  : B=4096,  F=1
#   block+3  : B=16384, F=1
#   block+4  : B=256,   F=4
#   block+5  : B=1024,  F=4
#   block+6  : B=4096,  F=4
#   block+7  : B=16384, F=4
#   block+8  : B=1024,  F=8
#   block+9  : B=4096,  F=8
#   block+10 : B=16384, F=8
$blockStarts = @(2, 13, 24, 35, 46)

# --- New column headers G1 / H1 ------------------------------------------
$ws.Range("G1").Value = "Speedup (1:4)"
$ws.Range("H1").Value = "Speedup (1:8)"

# --- "N/A" marker cells for the B=256 rows (no 8-thread data to compare) -
# These are filled in before the I1/J1 headers below so that the new
# shared-string table is built up in the same order the original author's
# edit produced: Speedup (1:4), Speedup (1:8), N/A, Efficency x2.
foreach ($start in $blockStarts) {
    $ws.Range("G$start").Value = "N/A"
    $ws.Range("H$start").Value = "N/A"
}

# --- New column headers I1 / J1 ------------------------------------------
$ws.Range("I1").Value = "Efficency (Speed: 4 Thread)"
$ws.Range("J1").Value = "Efficency (Speed: 8 Thread)"

# --- Speedup / Efficiency formulas ---------------------------------------
foreach ($start in $blockStarts) {
    $r2 = $start + 1    # B=1024,  F=1
    $r3 = $start + 2    # B=4096,  F=1
    $r4 = $start + 3    # B=16384, F=1

    $f4_1024  = $start + 5
    $f4_4096  = $start + 6
    $f4_16384 = $start + 7

    $f8_1024  = $start + 8
    $f8_4096  = $start + 9
    $f8_16384 = $start + 10

    $ws.Range("G$r2").Formula = "=E$r2/E$f4_1024"
    $ws.Range("H$r2").Formula = "=E$r2/E$f8_1024"
    $ws.Range("I$r2").Formula = "=G$r2/4"
    $ws.Range("J$r2").Formula = "=H$r2/8"

    $ws.Range("G$r3").Formula = "=E$r3/E$f4_4096"
    $ws.Range("H$r3").Formula = "=E$r3/E$f8_4096"
    $ws.Range("I$r3").Formula = "=G$r3/4"
    $ws.Range("J$r3").Formula = "=H$r3/8"

    $ws.Range("G$r4").Formula = "=E$r4/E$f4_16384"
    $ws.Range("H$r4").Formula = "=E$r4/E$f8_16384"
    $ws.Range("I$r4").Formula = "=G$r4/4"
    $ws.Range("J$r4").Formula = "=H$r4/8"
}

# --- New column widths ----------------------------------------------------
# (values chosen so the resulting stored width is the closest possible
# match to the author's final widths of 16.43 / 19 / 29.57 / 25.57 chars)
$ws.Columns.Item(7).ColumnWidth = 15.666666666666666
$ws.Columns.Item(8).ColumnWidth = 18.166666666666668
$ws.Columns.Item(9).ColumnWidth = 28.666666666666668
$ws.Columns.Item(10).ColumnWidth = 24.666666666666668

# --- Scroll / selection state (matches the author's final view) ---------
$ws.Application.Goto($ws.Range("A34"), $true)
$ws.Range("J53").Select()
